$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet (tab name + workbook.xml <sheet name="..."> entry)
$ws.Name = "o554F"

# 2) Append a new data row (row 16) holding the Gaussian-Quadrature-Scheme
#    averaged intensities (same layout/format as the existing HKL rows).
$row = 16
$prevRow = $row - 1

# Copy the formatting of the previous row's A cell (bold/centered/bordered
# "HKL index" style) onto the new row's A cell before setting its value.
$ws.Range("A$prevRow").Copy()
$ws.Range("A$row").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A$row").Value = 14
$ws.Range("B$row").Value = "HexGrid-60degTilt5degRes"

$values = @(
    1.07982118082457,
    0.7630299429227828,
    1.035803132436836,
    1.07982118082457,
    0.8744127469823707,
    1.100669706456135,
    1.050204621539189,
    0.7630299429227828,
    0.8994165376798096,
    0.9896188592521895,
    0.9839902218603139
)

$col = 3
foreach ($val in $values) {
    $ws.Cells.Item($row, $col).Value = $val
    $col++
}
